$d = $word.ActiveDocument

# Find every paragraph whose entire text is exactly "Phone - VARCHAR(20)"
# (the paragraph mark is included as Chr(13) at the end of Paragraph.Range.Text).
# This deliberately excludes the paragraph that already reads
# "Phone - VARCHAR(50) - Contact phone number." so only the three
# plain "Phone - VARCHAR(20)" bullets are touched, matching the diff.
$targetText = "Phone - VARCHAR(20)" + [char]13
$count = $d.Paragraphs.Count

for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -ne $targetText) {
        continue
    }

    $r = $p.Range
    $start = $r.Start

    # Layout of "Phone - VARCHAR(20)":
    #   [0,16)  -> "Phone - VARCHAR("   (kept, unchanged run)
    #   [16,17) -> "2"                  (becomes "5")
    #   [17,19) -> "0)"                 (kept, but split into its own run)

    # Temporarily give the leading "Phone - VARCHAR(" text different
    # character formatting than the rest of the run so that it will not
    # be recombined with its neighbours once their text changes.
    $prefixRange = $d.Range($start, $start + 16)
    $prefixRange.Font.Bold = 1

    # Replace the "2" with "5".
    $digitRange = $d.Range($start + 16, $start + 17)
    $digitRange.Text = "5"

    # Split "0)" away from "5" into its own run.
    $suffixRange = $d.Range($start + 17, $start + 19)
    $suffixRange.Font.Bold = 1
    $suffixRange.Font.Bold = 0

    # Restore the prefix's normal (non-bold) formatting.
    $prefixRange2 = $d.Range($start, $start + 16)
    $prefixRange2.Font.Bold = 0
}
